# Applies the value updates described by the commit diff to each leve-profit
# sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each touched row holds plain
# cached numbers (no formulas) in columns H..N, so the edits are direct
# Range.Value assignments; a handful of cells are removed/added entirely
# (ClearContents / new Value) to mirror cells disappearing or appearing in the
# underlying XML.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20458.396
$ws.Range("I32").Value = 3798.323
$ws.Range("J32").Value = 103758.766
$ws.Range("K32").Value = 3798.323
$ws.Range("L32").Value = 103758.766
$ws.Range("M32").Value = -3511.323
$ws.Range("N32").Value = -104332.766

$ws.Range("H44").Value = 37436.75
$ws.Range("J44").Value = 37436.75
$ws.Range("L44").Value = 37436.75
$ws.Range("N44").Value = -38412.75

$ws.Range("H55").Value = 22489.75
$ws.Range("J55").Value = 22489.75
$ws.Range("L55").Value = 22489.75
$ws.Range("N55").Value = -23119.75

$ws.Range("H80").Value = 36172
$ws.Range("J80").Value = 36172
$ws.Range("L80").Value = 36172
$ws.Range("N80").Value = -38168

$ws.Range("H83").Value = 36172
$ws.Range("J83").Value = 36172
$ws.Range("L83").Value = 108516
$ws.Range("N83").Value = -118500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18150.812
$ws.Range("J82").Value = 35450.855
$ws.Range("L82").Value = 35450.855
$ws.Range("N82").Value = -36216.855

$ws.Range("H85").Value = 18150.812
$ws.Range("J85").Value = 35450.855
$ws.Range("L85").Value = 35450.855
$ws.Range("N85").Value = -38102.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 27210
$ws.Range("J41").Value = 27210
$ws.Range("L41").Value = 27210
$ws.Range("N41").Value = -28066

$ws.Range("H50").Value = 21797.334
$ws.Range("J50").Value = 21797.334
$ws.Range("L50").Value = 21797.334
$ws.Range("N50").Value = -23047.334

$ws.Range("H60").Value = 20624.334
$ws.Range("J60").Value = 27936.5
$ws.Range("L60").Value = 27936.5
$ws.Range("N60").Value = -28958.5

$ws.Range("H68").Value = 25323.75
$ws.Range("J68").Value = 25323.75
$ws.Range("L68").Value = 25323.75
$ws.Range("N68").Value = -26821.75

$ws.Range("H71").Value = 25323.75
$ws.Range("J71").Value = 25323.75
$ws.Range("L71").Value = 75971.25
$ws.Range("N71").Value = -83459.25

$ws.Range("H134").Value = 10465.944
$ws.Range("I134").Value = 10465.944
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 31397.832
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -28862.832
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1539.5416
$ws.Range("I68").Value = 1083.591
$ws.Range("J68").Value = 1925.3462
$ws.Range("K68").Value = 3250.773
$ws.Range("L68").Value = 5776.0386
$ws.Range("M68").Value = -2439.773
$ws.Range("N68").Value = -7398.0386

$ws.Range("H71").Value = 1539.5416
$ws.Range("I71").Value = 1083.591
$ws.Range("J71").Value = 1925.3462
$ws.Range("K71").Value = 9752.319
$ws.Range("L71").Value = 17328.1158
$ws.Range("M71").Value = -5696.319
$ws.Range("N71").Value = -25440.1158

$ws.Range("H107").Value = 1414.3948
$ws.Range("I107").Value = 904.25
$ws.Range("J107").Value = 1649.8462
$ws.Range("K107").Value = 2712.75
$ws.Range("L107").Value = 4949.5386
$ws.Range("M107").Value = -792.75
$ws.Range("N107").Value = -8789.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H47").Value = 15333.333
$ws.Range("J47").Value = 15333.333
$ws.Range("L47").Value = 15333.333
$ws.Range("N47").Value = -16469.333

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H123").Value = 20881.666
$ws.Range("J123").Value = 20881.666
$ws.Range("L123").Value = 20881.666
$ws.Range("N123").Value = -25781.666

$ws.Range("H135").Value = 59342.855
$ws.Range("J135").Value = 59342.855
$ws.Range("L135").Value = 59342.855
$ws.Range("N135").Value = -69482.85500000001

$ws.Range("H138").Value = 39800
$ws.Range("J138").Value = 39800
$ws.Range("L138").Value = 39800
$ws.Range("N138").Value = -50080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5957.9165
$ws.Range("I7").Value = 2299.8
$ws.Range("J7").Value = 8570.857
$ws.Range("K7").Value = 2299.8
$ws.Range("L7").Value = 8570.857
$ws.Range("M7").Value = -2187.8
$ws.Range("N7").Value = -8794.857

$ws.Range("H22").Value = 855.65216
$ws.Range("I22").Value = 653.1667
$ws.Range("J22").Value = 1076.5454
$ws.Range("K22").Value = 653.1667
$ws.Range("L22").Value = 1076.5454
$ws.Range("M22").Value = -358.1667
$ws.Range("N22").Value = -1666.5454

$ws.Range("H27").Value = 855.65216
$ws.Range("I27").Value = 653.1667
$ws.Range("J27").Value = 1076.5454
$ws.Range("K27").Value = 653.1667
$ws.Range("L27").Value = 1076.5454
$ws.Range("M27").Value = -546.1667
$ws.Range("N27").Value = -1290.5454

$ws.Range("H126").Value = 5957.9165
$ws.Range("I126").Value = 2299.8
$ws.Range("J126").Value = 8570.857
$ws.Range("K126").Value = 6899.400000000001
$ws.Range("L126").Value = 25712.571
$ws.Range("M126").Value = -4429.400000000001
$ws.Range("N126").Value = -30652.571

$ws.Range("H132").Value = 1684.8206
$ws.Range("I132").Value = 1234.4
$ws.Range("J132").Value = 2489.1428
$ws.Range("K132").Value = 3703.2
$ws.Range("L132").Value = 7467.428400000001
$ws.Range("M132").Value = -1173.2
$ws.Range("N132").Value = -12527.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 17456.25
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 17456.25
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 17456.25
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -18212.25

$ws.Range("H51").Value = 8000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9020

$ws.Range("H57").Value = 15000
$ws.Range("J57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16508

$ws.Range("H58").Value = 10333.333
$ws.Range("I58").Value = 15000
$ws.Range("K58").Value = 15000
$ws.Range("M58").Value = -14692

Write-Output "Applied leve-profit updates across ARM/BSM/CRP/CUL/GSM/LTW/WVR."
